$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row so the old rows 7 & 8 shift down to 9 & 10 ---
$ws.Rows.Item(8).Insert()

# --- Update Profit (D) values for rows 2-6 ---
$ws.Cells.Item(2,4).Value = 488.74
$ws.Cells.Item(3,4).Value = 1113.05
$ws.Cells.Item(4,4).Value = -3.59
$ws.Cells.Item(5,4).Value = 559.9
$ws.Cells.Item(6,4).Value = 39.37

# --- Row 7 becomes a brand-new person (Nhung / NU / CSM - Ban hang) ---
$ws.Cells.Item(7,1).Value = "Nhung (NU)"
$ws.Cells.Item(7,2).Value = "NU"
$ws.Cells.Item(7,3).Value = "CSM - Bán hàng"
$ws.Cells.Item(7,4).Value = 5225.64
$ws.Cells.Item(7,5).Value = 5628.8

# --- Row 8 is the newly inserted person (Thien Ha / HV / CSM - Ban hang) ---
$ws.Cells.Item(8,1).Value = "Thiên Hà (HV)"
$ws.Cells.Item(8,2).Value = "HV"
$ws.Cells.Item(8,3).Value = "CSM - Bán hàng"
$ws.Cells.Item(8,4).Value = 5225.64
$ws.Cells.Item(8,5).Value = 5628.8

# --- Rows 9 & 10 keep the former rows 7 & 8 (Thanh / Truong), with new D values ---
$ws.Cells.Item(9,1).Value = "Thành  (BX)"
$ws.Cells.Item(9,2).Value = "BX"
$ws.Cells.Item(9,3).Value = "R&D"
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).Value = 1899.48

$ws.Cells.Item(10,1).Value = "Truong (XT)"
$ws.Cells.Item(10,2).Value = "XT"
$ws.Cells.Item(10,3).Value = "Designer"
$ws.Cells.Item(10,4).Value = 2748.07
$ws.Cells.Item(10,5).Value = 3769.29

# --- New column F: "KPI" header + percentage values stored as text ---
# (leading apostrophe forces the percentage strings to stay literal text,
#  matching the t="str" cells in the target workbook instead of being
#  auto-parsed into percent-formatted numbers)
$ws.Cells.Item(1,6).Value = "KPI"
$ws.Cells.Item(2,6).Value = "'34.49%"
$ws.Cells.Item(3,6).Value = "'34.07%"
$ws.Cells.Item(4,6).Value = "'-0.28%"
$ws.Cells.Item(5,6).Value = "'44.16%"
$ws.Cells.Item(6,6).Value = "'3.11%"
$ws.Cells.Item(7,6).Value = "'92.84%"
$ws.Cells.Item(8,6).Value = "'92.84%"
$ws.Cells.Item(9,6).Value = "'0.00%"
$ws.Cells.Item(10,6).Value = "'72.91%"
